# Update "想去人数" (F column) figures across the sheets, matching the
# regenerated site data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 92
$ws1.Range("F3").Value  = 1216
$ws1.Range("F4").Value  = 855
$ws1.Range("F5").Value  = 880
$ws1.Range("F6").Value  = 1607
$ws1.Range("F7").Value  = 338
$ws1.Range("F8").Value  = 1093
$ws1.Range("F11").Value = 224
$ws1.Range("F12").Value = 67
$ws1.Range("F13").Value = 566
$ws1.Range("F14").Value = 97
$ws1.Range("F15").Value = 59
$ws1.Range("F19").Value = 14
$ws1.Range("F20").Value = 605
$ws1.Range("F21").Value = 596
$ws1.Range("F22").Value = 84
$ws1.Range("F23").Value = 20
$ws1.Range("F24").Value = 804
$ws1.Range("F26").Value = 1
$ws1.Range("F27").Value = 214

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 213
$ws2.Range("F9").Value = 99

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 281

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 281
$ws4.Range("F4").Value  = 92
$ws4.Range("F5").Value  = 1216
$ws4.Range("F6").Value  = 855
$ws4.Range("F7").Value  = 880
$ws4.Range("F8").Value  = 1607
$ws4.Range("F9").Value  = 338
$ws4.Range("F10").Value = 1093
$ws4.Range("F13").Value = 224
$ws4.Range("F14").Value = 67
$ws4.Range("F15").Value = 566
$ws4.Range("F16").Value = 97
$ws4.Range("F17").Value = 59
$ws4.Range("F24").Value = 213
$ws4.Range("F25").Value = 213
$ws4.Range("F26").Value = 14
$ws4.Range("F27").Value = 605
$ws4.Range("F28").Value = 596
$ws4.Range("F29").Value = 84
$ws4.Range("F30").Value = 20
$ws4.Range("F31").Value = 804
$ws4.Range("F34").Value = 1
$ws4.Range("F35").Value = 214
$ws4.Range("F37").Value = 99
$ws4.Range("F38").Value = 99
